# Update the "Monitoramento diário" sheet, row 14 (day 21)
# with the reported daily counts: VMP=2, AMP=65, VMPP=65, AMPP=65
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitoramento diário")

$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 65
$ws.Range("G14").Value = 65
$ws.Range("I14").Value = 65

$ws.Activate()
[void]$ws.Range("I14").Select()
